# Updates the cryptos list (prices in column D, 1h volume % in column E)
# with freshly scraped values, and re-ranks Hedera/InternetComputer(DFINITY)
# which swapped places (rows 31 and 32).
#
# Values that look like plain numbers are written with a leading apostrophe
# so Excel keeps storing them as text (matching the original "text number"
# cells, e.g. preserving trailing zeros like "76.30") instead of silently
# converting them to numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.089.59'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '2.324.19'
$ws.Range('E3').Value = '  +2.75%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'253.87"
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('D6').Value = "'0.642"
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').Value = "'76.30"
$ws.Range('E7').Value = '  +6.85%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = "'0.656"
$ws.Range('E9').Value = '  -2.35%  '
$ws.Range('D10').Value = "'40.45"
$ws.Range('E10').Value = '  +1.59%  '
$ws.Range('E11').Value = '  +1.29%  '
$ws.Range('D12').Value = "'7.59"
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('D13').Value = "'0.106"
$ws.Range('E13').Value = '  +2.06%  '
$ws.Range('D14').Value = '2.671.52'
$ws.Range('E14').Value = '  +2.78%  '
$ws.Range('D15').Value = "'15.51"
$ws.Range('E15').Value = '  +4.21%  '
$ws.Range('D16').Value = "'0.885"
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').Value = '2.320.04'
$ws.Range('E17').Value = '  +2.63%  '
$ws.Range('D18').Value = '43.100.07'
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('E19').Value = '  +2.63%  '
$ws.Range('D20').Value = "'6.34"
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('D21').Value = "'73.12"
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('D22').Value = "'238.90"
$ws.Range('E22').Value = '  +0.36%  '
$ws.Range('D23').Value = "'2.24"
$ws.Range('E23').Value = '  +5.07%  '
$ws.Range('D24').Value = "'3.91"
$ws.Range('E24').Value = '  -0.98%  '
$ws.Range('D25').Value = "'11.65"
$ws.Range('E25').Value = '  -1.00%  '
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').Value = "'2.43"
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('E28').Value = '  -1.17%  '
$ws.Range('D29').Value = "'21.34"
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').Value = "'167.46"
$ws.Range('E30').Value = '  -0.39%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = "'0.0854"
$ws.Range('E31').Value = '  +10.42%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = "'6.32"
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('E33').Value = '  +0.76%  '
$ws.Range('D34').Value = "'30.63"
$ws.Range('E34').Value = '  +5.15%  '
$ws.Range('E35').Value = '  +1.84%  '
$ws.Range('D36').Value = "'4.60"
$ws.Range('E36').Value = '  +10.59%  '
$ws.Range('D37').Value = "'4.86"
$ws.Range('E37').Value = '  +2.33%  '
$ws.Range('E38').Value = '  -2.59%  '
$ws.Range('D39').Value = "'13.91"
$ws.Range('E39').Value = '  +13.60%  '
$ws.Range('D40').Value = "'2.35"
$ws.Range('E40').Value = '  +2.10%  '
$ws.Range('D41').Value = "'5.94"
$ws.Range('E41').Value = '  +1.49%  '
$ws.Range('E42').Value = '  +8.72%  '
$ws.Range('D43').Value = "'9.24"
$ws.Range('E43').Value = '  +3.28%  '
$ws.Range('D44').Value = "'62.82"
$ws.Range('E44').Value = '  -2.34%  '
$ws.Range('D45').Value = "'4.93"
$ws.Range('E45').Value = '  -2.12%  '
$ws.Range('D46').Value = "'106.13"
$ws.Range('E46').Value = '  +11.56%  '
$ws.Range('E47').Value = '  -0.50%  '
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('E49').Value = '  -0.23%  '
$ws.Range('D50').Value = "'1.19"
$ws.Range('E50').Value = '  -0.69%  '
$ws.Range('D51').Value = "'4.36"
$ws.Range('E51').Value = '  -0.46%  '
